# Adapt column header formatting to respective input file names:
#   *_old -> *_FV2210
#   *_new -> *_FV2304
# Then wrap the used range in a native Excel Table (adds xl/tables/table1.xml
# + tableParts on the sheet), and freeze the header row (pane split).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (row 1) cells: _old -> _FV2210, _new -> _FV2304 ---
$oldHeaders = @(
    "Segmentname_old",
    "Segmentgruppe_old",
    "Segment_old",
    "Datenelement_old",
    "Segment ID_old",
    "Code_old",
    "Qualifier_old",
    "Beschreibung_old",
    "Bedingungsausdruck_old",
    "Bedingung_old"
)
$newHeaders = @(
    "Segmentname_new",
    "Segmentgruppe_new",
    "Segment_new",
    "Datenelement_new",
    "Segment ID_new",
    "Code_new",
    "Qualifier_new",
    "Beschreibung_new",
    "Bedingungsausdruck_new",
    "Bedingung_new"
)

$fv2210Headers = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210"
)
$fv2304Headers = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

# Columns A-J (1-10) hold the "old" / FV2210 headers
for ($i = 0; $i -lt $fv2210Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2210Headers[$i]
}

# Column K (11) is "diff" - unchanged

# Columns L-U (12-21) hold the "new" / FV2304 headers
for ($i = 0; $i -lt $fv2304Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2304Headers[$i]
}

# --- 2. Turn the used range into an Excel Table ---
$rng = $ws.Range("A1:U71")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# --- 3. Freeze the header row (split pane under row 1) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
